# Applies the cryptos-list price/volume refresh described by the commit.
# Each text-like price/percentage cell is forced through a transient
# Text number format so Excel stores the exact literal string (matching
# values like "1.000" or "28.012.49" that are not valid round-trip
# numbers) instead of silently re-parsing it as a float. The format is
# reset to General immediately after so no stray cell styling remains.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

Set-TextValue "D2" "28.012.49"
Set-TextValue "E2" "  -0.62%  "

Set-TextValue "D3" "1.896.36"
Set-TextValue "E3" "  +1.40%  "

Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.02%  "

Set-TextValue "D5" "312.34"
Set-TextValue "E5" "  +0.17%  "

Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  -0.02%  "

Set-TextValue "D7" "0.5057"
Set-TextValue "E7" "  +0.37%  "

Set-TextValue "D8" "0.3921"
Set-TextValue "E8" "  -0.08%  "

Set-TextValue "D9" "0.09322"
Set-TextValue "E9" "  -2.82%  "

Set-TextValue "D10" "1.133"
Set-TextValue "E10" "  -0.46%  "

Set-TextValue "D11" "41.76"
Set-TextValue "E11" "  +2.19%  "

Set-TextValue "D12" "6.353"
Set-TextValue "E12" "  -2.23%  "

Set-TextValue "D13" "20.73"
Set-TextValue "E13" "  -1.31%  "

Set-TextValue "D14" "1.892.04"
Set-TextValue "E14" "  +1.14%  "

Set-TextValue "D15" "1.002"
Set-TextValue "E15" "  +0.03%  "

Set-TextValue "D16" "7.274"
Set-TextValue "E16" "  -1.97%  "

Set-TextValue "D17" "0.00001117"
Set-TextValue "E17" "  -1.18%  "

Set-TextValue "D18" "92.08"
Set-TextValue "E18" "  -1.03%  "

Set-TextValue "D19" "0.06578"
Set-TextValue "E19" "  -0.74%  "

Set-TextValue "D20" "17.74"
Set-TextValue "E20" "  +1.20%  "

Set-TextValue "E21" "  -0.11%  "

Set-TextValue "D22" "6.213"
Set-TextValue "E22" "  +0.90%  "

Set-TextValue "D23" "28.086.59"
Set-TextValue "E23" "  -0.59%  "

Set-TextValue "D24" "11.30"
Set-TextValue "E24" "  -0.20%  "

Set-TextValue "D25" "2.309"
Set-TextValue "E25" "  +1.17%  "

Set-TextValue "D26" "2.597"
Set-TextValue "E26" "  +2.69%  "

Set-TextValue "D27" "2.117.72"
Set-TextValue "E27" "  +1.55%  "

Set-TextValue "D28" "20.91"
Set-TextValue "E28" "  -1.33%  "

Set-TextValue "D29" "157.16"
Set-TextValue "E29" "  -0.34%  "

Set-TextValue "D30" "127.40"
Set-TextValue "E30" "  -0.15%  "

Set-TextValue "B31" "Stellar"
Set-TextValue "C31" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D31" "0.1068"
Set-TextValue "E31" "  +1.11%  "

Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "1.079"
Set-TextValue "E32" "  +1.00%  "

Set-TextValue "D33" "5.599"
Set-TextValue "E33" "  -0.46%  "

Set-TextValue "E34" "  -0.36%  "

Set-TextValue "D35" "9.560"
Set-TextValue "E35" "  -0.49%  "

Set-TextValue "D36" "0.06643"
Set-TextValue "E36" "  -1.62%  "

Set-TextValue "D37" "0.02405"
Set-TextValue "E37" "  +0.66%  "

Set-TextValue "B38" "Algorand"
Set-TextValue "C38" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2164"
Set-TextValue "E38" "  -0.99%  "

Set-TextValue "B39" "ARBITRUM"
Set-TextValue "C39" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D39" "1.219"
Set-TextValue "E39" "  -1.69%  "

Set-TextValue "D40" "1.266"
Set-TextValue "E40" "  +7.32%  "

Set-TextValue "D41" "0.6354"
Set-TextValue "E41" "  +0.08%  "

Set-TextValue "E42" "  +0.08%  "

Set-TextValue "D43" "11.39"
Set-TextValue "E43" "  -0.63%  "

Set-TextValue "E44" "  -0.11%  "

Set-TextValue "D45" "13.23"
Set-TextValue "E45" "  -2.84%  "

Set-TextValue "D46" "0.5971"
Set-TextValue "E46" "  -0.88%  "

Set-TextValue "D47" "3.706"
Set-TextValue "E47" "  +1.03%  "

Set-TextValue "E48" "  +0.52%  "

Set-TextValue "D49" "2.012"
Set-TextValue "E49" "  +1.14%  "

Set-TextValue "D50" "122.28"
Set-TextValue "E50" "  -1.43%  "

Set-TextValue "D51" "1.177"
Set-TextValue "E51" "  -1.57%  "
